$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; this shifts rows 20..57 down to 21..58,
# carrying along the existing cell styles (e.g. the date style on column D).
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with the new weekly price-report entry.
$ws.Cells.Item(20, 1).Value  = 11
$ws.Cells.Item(20, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value  = "Bíobío"
$ws.Cells.Item(20, 4).Value  = 45125
$ws.Cells.Item(20, 5).Value  = 8
$ws.Cells.Item(20, 6).Value  = 100114007
$ws.Cells.Item(20, 7).Value  = "Jengibre"
$ws.Cells.Item(20, 8).Value  = "Sin especificar"
$ws.Cells.Item(20, 9).Value  = "Primera"
$ws.Cells.Item(20, 10).Value = 30
$ws.Cells.Item(20, 11).Value = 16000
$ws.Cells.Item(20, 12).Value = 16000
$ws.Cells.Item(20, 13).Value = 16000
$ws.Cells.Item(20, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(20, 15).Value = "Perú"
$ws.Cells.Item(20, 16).Value = 1231
$ws.Cells.Item(20, 17).Value = 13
$ws.Cells.Item(20, 18).Value = "Hortaliza"
